# Update training hyper-parameters (Train25 -> Train26) on the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# batch_size: 128 -> 32
$ws.Range("B2").Value = 32

# lr: 0.01 -> 0.0002
$ws.Range("B9").Value = 0.0002

# beta1: 0.5 -> 0.2
$ws.Range("B10").Value = 0.2

# n_epochs: 1 -> 100
$ws.Range("B13").Value = 100
